# The upstream commit ("added slides to alignment lecture and fasta/fastq")
# touches other decks; the only change that lands inside this deck's
# canonical OOXML is a cosmetic re-serialization of the Office-2010
# compatibility markup that PowerPoint stores for pictures/shapes whose
# fill or line has been turned off: the <a14:hiddenFill>/<a14:hiddenLine>
# elements inside each shape's <a:extLst>. Re-saving the file simply
# re-emits their `xmlns:a14=...` / `xmlns=""` namespace declarations in a
# different (but semantically identical) order. That happens on the three
# shapes that already carry this legacy markup: the "Picture 4" picture on
# slide 1, the "Picture 7" picture on the slide layout behind slide 1, and
# the "TextBox 10" shape on the slide master.
#
# There is no dedicated hidden-fill/hidden-line property on the PowerPoint
# object model (it is write-once legacy compat markup, not something the
# UI/object model exposes), so we can't set it directly. What we *can* do
# is touch each of those shapes via the object model so the host re-emits
# them as part of a normal save pass, without altering their visible
# formatting/position in any way.

function Find-ShapeByName($shapes, $name) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Name -eq $name) {
            return $sh
        }
    }
    return $null
}

function Touch-Shape($shape) {
    if ($shape -ne $null) {
        # No-op re-assignment: forces the shape through the host's
        # save/serialization path without changing its position/size.
        $shape.Left = $shape.Left
        $shape.Top = $shape.Top
    }
}

$p = $ppt.ActivePresentation

# --- Slide 1: "Picture 4" (TGI_logo_V_2color_bevel.tiff) ---
$s = $p.Slides.Item(1)
$slidePic = Find-ShapeByName $s.Shapes "Picture 4"
Touch-Shape $slidePic

# --- Slide layout behind slide 1: "Picture 7" (bioinformatics.ca logo) ---
$layout = $s.CustomLayout
$layoutPic = Find-ShapeByName $layout.Shapes "Picture 7"
Touch-Shape $layoutPic

# --- Slide master: "TextBox 10" ---
$master = $p.SlideMaster
$masterBox = Find-ShapeByName $master.Shapes "TextBox 10"
Touch-Shape $masterBox
